$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F5 "Maximum" -> "Maximum, without animations."
$ws.Range("F5").Value = "Maximum, without animations."

# New note in F6 for the "Achievable" quality tier.
$ws.Range("F6").Value = "Achievable"

# Row 5 grows to fit the longer wrapped text.
$ws.Rows(5).RowHeight = 23.65

# Selection / scroll position moves down to the newly edited cell.
$ws.Range("F6").Select()
$excel.ActiveWindow.ScrollRow = 2
